$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 14006.5
$ws.Range("J17").Value = 21590.4
$ws.Range("L17").Value = 64771.2
$ws.Range("N17").Value = -65107.2

$ws.Range("H32").Value = 5206.769
$ws.Range("I32").Value = 11005
$ws.Range("K32").Value = 11005
$ws.Range("M32").Value = -10679

$ws.Range("H53").Value = 1402.0555
$ws.Range("I53").Value = 278.16666
$ws.Range("J53").Value = 1964
$ws.Range("K53").Value = 278.16666
$ws.Range("L53").Value = 1964
$ws.Range("M53").Value = 358.83334
$ws.Range("N53").Value = -3238

$ws.Range("H101").Value = 2064.3333
$ws.Range("J101").Value = 2087.5
$ws.Range("L101").Value = 6262.5
$ws.Range("N101").Value = -9506.5

$ws.Range("H112").Value = 1224.7778
$ws.Range("I112").Value = 766.8570999999999
$ws.Range("J112").Value = 1335.3103
$ws.Range("K112").Value = 2300.5713
$ws.Range("L112").Value = 4005.9309
$ws.Range("M112").Value = -1192.5713
$ws.Range("N112").Value = -6221.9309

$ws.Range("H116").Value = 6653
$ws.Range("I116").Value = 2794.6667
$ws.Range("K116").Value = 2794.6667
$ws.Range("M116").Value = 647.3332999999998

$ws.Range("H132").Value = 2117.0952
$ws.Range("I132").Value = 2047.6571
$ws.Range("J132").Value = 2464.2856
$ws.Range("K132").Value = 6142.971299999999
$ws.Range("L132").Value = 7392.8568
$ws.Range("M132").Value = -3612.971299999999
$ws.Range("N132").Value = -12452.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8886.762000000001
$ws.Range("I32").Value = 7176.1333
$ws.Range("K32").Value = 7176.1333
$ws.Range("M32").Value = -6889.1333

$ws.Range("H92").Value = 34275
$ws.Range("J92").Value = 34275
$ws.Range("L92").Value = 34275
$ws.Range("N92").Value = -39267

$ws.Range("H97").Value = 1100.5834
$ws.Range("I97").Value = 800.8889
$ws.Range("J97").Value = 1999.6666
$ws.Range("K97").Value = 800.8889
$ws.Range("L97").Value = 1999.6666
$ws.Range("M97").Value = -304.8889
$ws.Range("N97").Value = -2991.6666

$ws.Range("H110").Value = 4633.6
$ws.Range("I110").Value = 4633.6
$ws.Range("K110").Value = 4633.6
$ws.Range("M110").Value = -2588.6

$ws.Range("H131").Value = 46278.2
$ws.Range("J131").Value = 46278.2
$ws.Range("L131").Value = 46278.2
$ws.Range("N131").Value = -56358.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3077.524
$ws.Range("I86").Value = 2614.2666
$ws.Range("K86").Value = 2614.2666
$ws.Range("M86").Value = -1491.2666

$ws.Range("H89").Value = 3077.524
$ws.Range("I89").Value = 2614.2666
$ws.Range("K89").Value = 13071.333
$ws.Range("M89").Value = -7455.332999999999

$ws.Range("H105").Value = 33054.668
$ws.Range("I105").Value = 52248.25
$ws.Range("J105").Value = 17699.8
$ws.Range("K105").Value = 52248.25
$ws.Range("L105").Value = 17699.8
$ws.Range("M105").Value = -50501.25
$ws.Range("N105").Value = -21193.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 521.6
$ws.Range("I2").Value = 1578
$ws.Range("J2").Value = 257.5
$ws.Range("K2").Value = 1578
$ws.Range("L2").Value = 257.5
$ws.Range("M2").Value = -1465
$ws.Range("N2").Value = -483.5

$ws.Range("H8").Value = 936.3333
$ws.Range("I8").Value = 9
$ws.Range("J8").Value = 1400
$ws.Range("K8").Value = 9
$ws.Range("L8").Value = 1400
$ws.Range("M8").Value = 131
$ws.Range("N8").Value = -1680

$ws.Range("H10").Value = 469.75
$ws.Range("I10").Value = 185
$ws.Range("K10").Value = 185
$ws.Range("M10").Value = -46

$ws.Range("H11").Value = 620.25
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 620.25
$ws.Range("K11").Value = 0
$ws.Range("M11").Value = 620.25
$ws.Range("N11").Value = -900.25

$ws.Range("H13").Value = 2505
$ws.Range("J13").Value = 2505
$ws.Range("L13").Value = 2505
$ws.Range("N13").Value = -2783

$ws.Range("H22").Value = 458.54544
$ws.Range("I22").Value = 317.66666
$ws.Range("J22").Value = 1092.5
$ws.Range("K22").Value = 317.66666
$ws.Range("L22").Value = 1092.5
$ws.Range("M22").Value = 32.33334000000002
$ws.Range("N22").Value = -1792.5

$ws.Range("H31").Value = 44958.07
$ws.Range("I31").Value = 6803.6924
$ws.Range("J31").Value = 78025.2
$ws.Range("K31").Value = 6803.6924
$ws.Range("L31").Value = 78025.2
$ws.Range("M31").Value = -6508.6924
$ws.Range("N31").Value = -78615.2

$ws.Range("H34").Value = 44958.07
$ws.Range("I34").Value = 6803.6924
$ws.Range("J34").Value = 78025.2
$ws.Range("K34").Value = 6803.6924
$ws.Range("L34").Value = 78025.2
$ws.Range("M34").Value = -6601.6924
$ws.Range("N34").Value = -78429.2

$ws.Range("H58").Value = 5818
$ws.Range("I58").Value = 2311.889
$ws.Range("J58").Value = 16336.333
$ws.Range("K58").Value = 2311.889
$ws.Range("L58").Value = 16336.333
$ws.Range("M58").Value = -2108.889
$ws.Range("N58").Value = -16742.333

$ws.Range("H99").Value = 3750.889
$ws.Range("I99").Value = 3498.75
$ws.Range("J99").Value = 3952.6
$ws.Range("K99").Value = 3498.75
$ws.Range("L99").Value = 3952.6
$ws.Range("M99").Value = -2000.75
$ws.Range("N99").Value = -6948.6

$ws.Range("H107").Value = 884.79486
$ws.Range("J107").Value = 1095.8334
$ws.Range("L107").Value = 1095.8334
$ws.Range("N107").Value = -4935.8334

$ws.Range("H126").Value = 3750.889
$ws.Range("I126").Value = 3498.75
$ws.Range("J126").Value = 3952.6
$ws.Range("K126").Value = 10496.25
$ws.Range("L126").Value = 11857.8
$ws.Range("M126").Value = -8026.25
$ws.Range("N126").Value = -16797.8

$ws.Range("H134").Value = 2752.7646
$ws.Range("I134").Value = 2143.125
$ws.Range("K134").Value = 6429.375
$ws.Range("M134").Value = -3894.375

$ws.Range("H136").Value = 5818
$ws.Range("I136").Value = 2311.889
$ws.Range("J136").Value = 16336.333
$ws.Range("K136").Value = 6935.667
$ws.Range("L136").Value = 49008.999
$ws.Range("M136").Value = -4385.667
$ws.Range("N136").Value = -54108.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1640.1765
$ws.Range("J113").Value = 1676.1666
$ws.Range("L113").Value = 5028.4998
$ws.Range("N113").Value = -9368.4998

$ws.Range("H122").Value = 1657.64
$ws.Range("I122").Value = 481.15384
$ws.Range("J122").Value = 2932.1667
$ws.Range("K122").Value = 4330.38456
$ws.Range("L122").Value = 26389.5003
$ws.Range("M122").Value = -1880.38456
$ws.Range("N122").Value = -31289.5003

$ws.Range("H131").Value = 7430750.5
$ws.Range("J131").Value = 5257153
$ws.Range("L131").Value = 15771459
$ws.Range("N131").Value = -15781539

$ws.Range("H136").Value = 2853.25
$ws.Range("I136").Value = 2658.0908
$ws.Range("K136").Value = 7974.2724
$ws.Range("M136").Value = -2874.2724

$ws.Range("H137").Value = 85499.086
$ws.Range("I137").Value = 2200
$ws.Range("J137").Value = 113265.445
$ws.Range("K137").Value = 6600
$ws.Range("L137").Value = 339796.335
$ws.Range("M137").Value = -1500
$ws.Range("N137").Value = -349996.335

$ws.Range("H139").Value = 3005.257
$ws.Range("I139").Value = 2190.2354
$ws.Range("J139").Value = 3775
$ws.Range("K139").Value = 6570.706200000001
$ws.Range("L139").Value = 11325
$ws.Range("M139").Value = -1430.706200000001
$ws.Range("N139").Value = -21605

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13639.692
$ws.Range("I70").Value = 10543.643
$ws.Range("J70").Value = 17251.75
$ws.Range("K70").Value = 10543.643
$ws.Range("L70").Value = 17251.75
$ws.Range("M70").Value = -10273.643
$ws.Range("N70").Value = -17791.75

$ws.Range("H73").Value = 13639.692
$ws.Range("I73").Value = 10543.643
$ws.Range("J73").Value = 17251.75
$ws.Range("K73").Value = 10543.643
$ws.Range("L73").Value = 17251.75
$ws.Range("M73").Value = -9607.643
$ws.Range("N73").Value = -19123.75

$ws.Range("H80").Value = 359989.5
$ws.Range("I80").Value = 502476.4
$ws.Range("K80").Value = 502476.4
$ws.Range("M80").Value = -501478.4

$ws.Range("H83").Value = 359989.5
$ws.Range("I83").Value = 502476.4
$ws.Range("K83").Value = 2512382
$ws.Range("M83").Value = -2507390

$ws.Range("H122").Value = 6086.515
$ws.Range("I122").Value = 5709.174
$ws.Range("J122").Value = 6954.4
$ws.Range("K122").Value = 17127.522
$ws.Range("L122").Value = 20863.2
$ws.Range("M122").Value = -14677.522
$ws.Range("N122").Value = -25763.2

$ws.Range("H126").Value = 4344.4443
$ws.Range("I126").Value = 3300
$ws.Range("K126").Value = 9900
$ws.Range("M126").Value = -7430

$ws.Range("H132").Value = 8581.777
$ws.Range("I132").Value = 7487.2915
$ws.Range("J132").Value = 17337.666
$ws.Range("K132").Value = 22461.8745
$ws.Range("L132").Value = 52012.99800000001
$ws.Range("M132").Value = -19931.8745
$ws.Range("N132").Value = -57072.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4201.88
$ws.Range("I22").Value = 2279
$ws.Range("J22").Value = 7620.3335
$ws.Range("K22").Value = 2279
$ws.Range("L22").Value = 7620.3335
$ws.Range("M22").Value = -1984
$ws.Range("N22").Value = -8210.333500000001

$ws.Range("H27").Value = 4201.88
$ws.Range("I27").Value = 2279
$ws.Range("J27").Value = 7620.3335
$ws.Range("K27").Value = 2279
$ws.Range("L27").Value = 7620.3335
$ws.Range("M27").Value = -2172
$ws.Range("N27").Value = -7834.3335

$ws.Range("H122").Value = 163643.36
$ws.Range("I122").Value = 213035
$ws.Range("K122").Value = 639105
$ws.Range("M122").Value = -636655

$ws.Range("H131").Value = 46666
$ws.Range("I131").Value = 39998
$ws.Range("J131").Value = 50000
$ws.Range("K131").Value = 39998
$ws.Range("L131").Value = 50000
$ws.Range("M131").Value = -34958
$ws.Range("N131").Value = -60080

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2436.1875
$ws.Range("I122").Value = 2274.625
$ws.Range("K122").Value = 6823.875
$ws.Range("M122").Value = -4373.875

$ws.Range("H136").Value = 4849.4546
$ws.Range("I136").Value = 3433.9
$ws.Range("J136").Value = 19005
$ws.Range("K136").Value = 10301.7
$ws.Range("L136").Value = 57015
$ws.Range("M136").Value = -7751.700000000001
$ws.Range("N136").Value = -62115
